$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 286
$ws.Range("Q3").Value = 2391

$ws.Range("M4").Value = 80
$ws.Range("Q4").Value = 445

$ws.Range("M5").Value = 366
$ws.Range("Q5").Value = 2836

$ws.Range("M8").Value = 53
$ws.Range("Q8").Value = 334

$ws.Range("M9").Value = 7
$ws.Range("Q9").Value = 31

$ws.Range("M10").Value = 16
$ws.Range("Q10").Value = 71

$ws.Range("M13").Value = 20
$ws.Range("Q13").Value = 132

$ws.Range("M14").Value = 96
$ws.Range("Q14").Value = 571

$ws.Range("M18").Value = 37
$ws.Range("Q18").Value = 169

$ws.Range("M20").Value = 2
$ws.Range("Q20").Value = 53

$ws.Range("M21").Value = 75
$ws.Range("Q21").Value = 411

$ws.Range("M22").Value = 13
$ws.Range("Q22").Value = 26

$ws.Range("M23").Value = 2
$ws.Range("Q23").Value = 11

$ws.Range("M24").Value = 3
$ws.Range("Q24").Value = 19

$ws.Range("M25").Value = 1
$ws.Range("Q25").Value = 2

$ws.Range("M26").Value = 5
$ws.Range("Q26").Value = 19

$ws.Range("M28").Value = 12
$ws.Range("Q28").Value = 115

$ws.Range("M29").Value = 26
$ws.Range("Q29").Value = 110
